$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ARDEL Deputy Integration Lead"
$ws.Range("A3").Value = "AMDR Test Director"
$ws.Range("A4").Value = "ARDEL Integration Lead / Site Lead"
$ws.Range("A5").Value = "CSEDS Integration Lead"
$ws.Range("A6").Value = "AMDR Integration Lead"

$ws.Columns("A:A").AutoFit()

$ws.Range("A15").Select()
